$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.557.16'
$ws.Range("E2").Value = '  +1.96%  '

$ws.Range("D3").Value = '1.666.17'
$ws.Range("E3").Value = '  +0.86%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9987'
$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.01'
$ws.Range("E5").Value = '  -0.12%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4805'
$ws.Range("E7").Value = '  -0.26%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2635'
$ws.Range("E8").Value = '  +0.34%  '

$ws.Range("E9").Value = '  +2.51%  '

$ws.Range("E10").Value = '  -1.16%  '

$ws.Range("D11").Value = '1.660.83'
$ws.Range("E11").Value = '  +0.58%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.83'
$ws.Range("E12").Value = '  -0.69%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5903'
$ws.Range("E13").Value = '  -5.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.380'
$ws.Range("E14").Value = '  -5.01%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '75.22'
$ws.Range("E15").Value = '  +2.86%  '

$ws.Range("E16").Value = '  +0.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9993'
$ws.Range("E17").Value = '  -0.01%  '

$ws.Range("D18").Value = '25.540.21'
$ws.Range("E18").Value = '  +1.93%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006760'
$ws.Range("E19").Value = '  +2.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.49'
$ws.Range("E20").Value = '  -0.61%  '

$ws.Range("D21").Value = '1.874.18'
$ws.Range("E21").Value = '  +0.88%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.430'
$ws.Range("E22").Value = '  -2.49%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.722'
$ws.Range("E23").Value = '  +0.98%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.297'
$ws.Range("E24").Value = '  -0.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.98'
$ws.Range("E25").Value = '  +2.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.08'
$ws.Range("E26").Value = '  +0.29%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.380'
$ws.Range("E27").Value = '  -1.17%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '105.39'
$ws.Range("E28").Value = '  +1.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.712'
$ws.Range("E29").Value = '  +1.58%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.954'
$ws.Range("E30").Value = '  +4.31%  '

$ws.Range("E31").Value = '  +1.47%  '

$ws.Range("E32").Value = '  -2.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9991'

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04239'
$ws.Range("E34").Value = '  -8.29%  '

$ws.Range("E35").Value = '  +0.27%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6126'
$ws.Range("E36").Value = '  +5.76%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9519'
$ws.Range("E37").Value = '  +0.81%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.601'
$ws.Range("E38").Value = '  -0.23%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.8618'
$ws.Range("E39").Value = '  +3.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9993'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.858'
$ws.Range("E41").Value = '  +1.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01468'
$ws.Range("E42").Value = '  -6.31%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '97.09'
$ws.Range("E43").Value = '  -1.39%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3775'
$ws.Range("E44").Value = '  +0.68%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.849'
$ws.Range("E45").Value = '  +0.95%  '

$ws.Range("E46").Value = '  -1.77%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.218'
$ws.Range("E47").Value = '  +1.26%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05261'
$ws.Range("E48").Value = '  +1.38%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.78'
$ws.Range("E49").Value = '  -0.23%  '

$ws.Range("B50").Value = 'TrueUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.002'
$ws.Range("E50").Value = '  +0.06%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.339'
$ws.Range("E51").Value = '  +1.27%  '
